$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "Date" values (column D) for rows 2 through 33 (the header in
# D1 stays untouched), leaving the cells present but empty.
$ws.Range("D2:D33").ClearContents()
